# Two new data rows were recorded for "Paine" (Zapallo, Feria Lagunitas de Puerto
# Montt) with date 44596 (2022-02-08). They were inserted into the data table
# right before what was row 157, pushing all the existing rows (old 157..253)
# down by two (to 159..255).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 157 (shifts old row 157 -> 159, etc.)
$ws.Rows.Item(157).Insert()
$ws.Rows.Item(157).Insert()

# Populate the first new row (157)
$ws.Range("A157").Value() = 4
$ws.Range("B157").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C157").Value() = "Los Lagos"
$ws.Range("D157").Value() = 44596
$ws.Range("E157").Value() = 10
$ws.Range("F157").Value() = 100112045
$ws.Range("G157").Value() = "Zapallo"
$ws.Range("H157").Value() = "Paine"
$ws.Range("I157").Value() = "1a nueva(o)"
$ws.Range("J157").Value() = 600
$ws.Range("K157").Value() = 500
$ws.Range("L157").Value() = 500
$ws.Range("M157").Value() = 500
$ws.Range("N157").Value() = "$/kilo (volumen en unidades)"
$ws.Range("O157").Value() = "Región de O'Higgins"
$ws.Range("P157").Value() = 500
$ws.Range("Q157").Value() = 1
$ws.Range("R157").Value() = "Hortaliza"

# Populate the second new row (158)
$ws.Range("A158").Value() = 4
$ws.Range("B158").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C158").Value() = "Los Lagos"
$ws.Range("D158").Value() = 44596
$ws.Range("E158").Value() = 10
$ws.Range("F158").Value() = 100112045
$ws.Range("G158").Value() = "Zapallo"
$ws.Range("H158").Value() = "Paine"
$ws.Range("I158").Value() = "2a nueva(o)"
$ws.Range("J158").Value() = 600
$ws.Range("K158").Value() = 400
$ws.Range("L158").Value() = 400
$ws.Range("M158").Value() = 400
$ws.Range("N158").Value() = "$/kilo (volumen en unidades)"
$ws.Range("O158").Value() = "Región de O'Higgins"
$ws.Range("P158").Value() = 400
$ws.Range("Q158").Value() = 1
$ws.Range("R158").Value() = "Hortaliza"

Write-Host "Rows inserted and populated"
